# update get service id (dien+nuoc)
# Renames the "Apartment" id codes in column A to a simple sequential
# A001..A020 scheme and extends the meter-reading list from 9 rows
# (rows 2-10) to 20 rows (rows 2-21), reusing the existing
# usage(B)/unit(C) pattern for the newly added rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: renumber existing rows 2-10 to A001..A009 -------------
$ws.Range("A2").Value  = "A001"
$ws.Range("A3").Value  = "A002"
$ws.Range("A4").Value  = "A003"
$ws.Range("A5").Value  = "A004"
$ws.Range("A6").Value  = "A005"
$ws.Range("A7").Value  = "A006"
$ws.Range("A8").Value  = "A007"
$ws.Range("A9").Value  = "A008"
$ws.Range("A10").Value = "A009"

# Rows 4-10 pick up the same cell style already used by rows 2-3.
$ws.Range("A4:A10").Style = "Hyperlink"

# --- Add the new rows 11-21 (A010..A020) ------------------------------
$newRows = @(
    @{ Row = 11; Id = "A010"; Usage = 33 },
    @{ Row = 12; Id = "A011"; Usage = 12 },
    @{ Row = 13; Id = "A012"; Usage = 43 },
    @{ Row = 14; Id = "A013"; Usage = 21 },
    @{ Row = 15; Id = "A014"; Usage = 22 },
    @{ Row = 16; Id = "A015"; Usage = 23 },
    @{ Row = 17; Id = "A016"; Usage = 42 },
    @{ Row = 18; Id = "A017"; Usage = 14 },
    @{ Row = 19; Id = "A018"; Usage = 11 },
    @{ Row = 20; Id = "A019"; Usage = 14 },
    @{ Row = 21; Id = "A020"; Usage = 11 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Id
    $ws.Cells.Item($row, 2).Value = $r.Usage
    $ws.Cells.Item($row, 3).Value = "m^3"

    $ws.Cells.Item($row, 1).Style = "Hyperlink"
    $ws.Range($ws.Cells.Item($row, 2), $ws.Cells.Item($row, 3)).Style = "Normal"
}

# --- View state: selection now sits on the last added row ------------
$ws.Range("C20:C21").Select()
